$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.665.87'
$ws.Range("E2").Value = '  -0.69%  '

# Row 3
$ws.Range("D3").Value = '3.794.67'
$ws.Range("E3").Value = '  -2.01%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '700.68'
$ws.Range("E5").Value = '  +0.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.10'
$ws.Range("E6").Value = '  -2.84%  '

# Row 7
$ws.Range("D7").Value = '3.792.37'
$ws.Range("E7").Value = '  -2.04%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  -1.21%  '

# Row 10
$ws.Range("E10").Value = '  -2.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.56'
$ws.Range("E11").Value = '  +5.83%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  -0.96%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("E13").Value = '  -3.43%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.62'
$ws.Range("E14").Value = '  -2.44%  '

# Row 15
$ws.Range("D15").Value = '4.435.31'
$ws.Range("E15").Value = '  -1.98%  '

# Row 16
$ws.Range("D16").Value = '3.773.72'
$ws.Range("E16").Value = '  -2.46%  '

# Row 17
$ws.Range("D17").Value = '70.683.26'
$ws.Range("E17").Value = '  -0.81%  '

# Row 18
$ws.Range("E18").Value = '  +0.37%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.29'
$ws.Range("E19").Value = '  -2.40%  '

# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.08'
$ws.Range("E20").Value = '  -2.36%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '492.90'
$ws.Range("E21").Value = '  -1.53%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.64'
$ws.Range("E22").Value = '  -4.92%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("E23").Value = '  -0.64%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.11'
$ws.Range("E24").Value = '  -0.99%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000142'
$ws.Range("E25").Value = '  -4.75%  '

# Row 26
$ws.Range("D26").Value = '3.948.05'
$ws.Range("E26").Value = '  -1.63%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.00'
$ws.Range("E27").Value = '  -2.34%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.23'
$ws.Range("E28").Value = '  -5.78%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.01'
$ws.Range("E30").Value = '  -7.02%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.00'
$ws.Range("E31").Value = '  -6.10%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("E32").Value = '  -0.77%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.28'
$ws.Range("E33").Value = '  -4.46%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.90'
$ws.Range("E34").Value = '  -2.87%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.175'
$ws.Range("E35").Value = '  -3.17%  '

# Row 36
$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.09%  '

# Row 37
$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.768.77'
$ws.Range("E37").Value = '  -1.45%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.02'
$ws.Range("E38").Value = '  -2.63%  '

# Row 39
$ws.Range("E39").Value = '  -3.95%  '

# Row 40
$ws.Range("E40").Value = '  -2.83%  '

# Row 41
$ws.Range("E41").Value = '  -2.88%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.91'
$ws.Range("E42").Value = '  -1.93%  '

# Row 43
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.02%  '

# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.21'
$ws.Range("E44").Value = '  -6.51%  '

# Row 45
$ws.Range("E45").Value = '  +0.12%  '

# Row 46
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '166.75'
$ws.Range("E46").Value = '  +1.78%  '

# Row 47
$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000315'
$ws.Range("E47").Value = '  +0.97%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.90'
$ws.Range("E48").Value = '  +0.03%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '421.00'
$ws.Range("E49").Value = '  +0.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.55'
$ws.Range("E50").Value = '  -1.51%  '

# Row 51
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.291'
$ws.Range("E51").Value = '  -4.29%  '
